# Add data for 2021-10-11 (update "through" date in sheet name and
# October row label, plus refreshed counts/rates for row 12 (October
# partial month) and row 13 (Total) on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2021-10-03"

# Row 12: October (through 10-03) partial-month figures.
$ws.Range("A12").Value = "October (through 10-03)"
$ws.Range("I12").Value = 8
$ws.Range("J12").Value = 0.2
$ws.Range("L12").Value = 8
$ws.Range("O12").Value = 2
$ws.Range("R12").Value = 13
$ws.Range("U12").Value = 30

# Row 13: Total row, updated to include the new data.
$ws.Range("I13").Value = 585
$ws.Range("J13").Value = 0.0816
$ws.Range("L13").Value = 495
$ws.Range("M13").Value = 0.1097
$ws.Range("O13").Value = 381
$ws.Range("P13").Value = 0.1014
$ws.Range("R13").Value = 861
$ws.Range("S13").Value = 0.058
$ws.Range("U13").Value = 1200
$ws.Range("V13").Value = 0.061
